{"js": "// 1) Text edit: insert \"/**+-*+-*\" right after \"the\" in\n//    \"...How does my class compare to the rest of the country?\"\nconst body = context.document.body;\n\nconst searchResults = body.search(\"the rest of the country?\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"the rest of the/**+-*+-* country?\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 2) Remove the two trailing near-empty paragraphs (a non-breaking-space-only\n//    paragraph and the blank paragraph right after it) that sit right before\n//    \"Declare the statistical question you will report on here:\" - while\n//    leaving the preceding (identical-looking) blank/space paragraphs intact.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Declare the statistical question\") === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  const pBlank = target.getPreviousOrNullObject();\n  await context.sync();\n\n  const pSpace = pBlank.getPreviousOrNullObject();\n  await context.sync();\n\n  pBlank.delete();\n  pSpace.delete();\n  await context.sync();\n}\n", "ps1": "# 1) Text edit: insert \"/**+-*+-*\" right after \"the\" in\n#    \"...How does my class compare to the rest of the country?\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"the rest of the country?\"\n$find.Replacement.Text = \"the rest of the/**+-*+-* country?\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\n# 2) Remove the two trailing near-empty paragraphs (a non-breaking-space-only\n#    paragraph and the blank paragraph right after it) that sit right before\n#    \"Declare the statistical question you will report on here:\" - while\n#    leaving the preceding (identical-looking) blank/space paragraphs intact.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Declare the statistical question*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $pBlank = $target.Previous(1)\n    $pSpace = $pBlank.Previous(1)\n\n    $delRange = $d.Range($pSpace.Range.Start, $pBlank.Range.End)\n    $delRange.Delete()\n}\n"}
